# Applies the "data up to 22" update to the fb-surveyState workbook.
# - Corrects a handful of recomputed percentage values in rows 132-136
# - Backfills the AR column (Puerto Rico) for rows 114-119
# - Adds survey data for 17-22 June 2020 (rows 137-144)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Recalculated value corrections (rows 132-136) ---
$ws.Range("AM132").Value = 0.3668641
$ws.Range("AM133").Value = 0.341182
$ws.Range("AM134").Value = 0.3360986
$ws.Range("H135").Value = 0.323868
$ws.Range("AM135").Value = 0.3423892
$ws.Range("C136").Value = 0.5914961
$ws.Range("G136").Value = 0.357074
$ws.Range("H136").Value = 0.3257062
$ws.Range("M136").Value = 0.4520485
$ws.Range("U136").Value = 0.4387699
$ws.Range("AM136").Value = 0.3721072
$ws.Range("AP136").Value = 0.3558797

# --- Backfill Puerto Rico (AR) for rows 114-119 ---
$ws.Range("AR114").Value = 0
$ws.Range("AR115").Value = 0
$ws.Range("AR116").Value = 0.9259259
$ws.Range("AR117").Value = 0.9259259
$ws.Range("AR118").Value = 0.952381
$ws.Range("AR119").Value = 0.9090909

# --- New survey rows 137-144 (17-22 June 2020) ---
$ws.Range("A137").Value = "15 06 2020"
$ws.Range("B137").Value = 0.3863543
$ws.Range("C137").Value = 0.5801948
$ws.Range("D137").Value = 0.4773163
$ws.Range("F137").Value = 0.6350183
$ws.Range("G137").Value = 0.3609408
$ws.Range("H137").Value = 0.3065941
$ws.Range("I137").Value = 0.2844822
$ws.Range("J137").Value = 0.2401838
$ws.Range("K137").Value = 0.1973559
$ws.Range("L137").Value = 0.3593247
$ws.Range("M137").Value = 0.4467009
$ws.Range("O137").Value = 0.2677376
$ws.Range("P137").Value = 0.4920518
$ws.Range("Q137").Value = 0.4505365
$ws.Range("R137").Value = 0.3289479
$ws.Range("S137").Value = 0.5274879
$ws.Range("T137").Value = 0.3879951
$ws.Range("U137").Value = 0.4598094
$ws.Range("V137").Value = 0.5519435
$ws.Range("W137").Value = 0.2722746
$ws.Range("X137").Value = 0.3334819
$ws.Range("Y137").Value = 0.3761819
$ws.Range("Z137").Value = 0.268493
$ws.Range("AA137").Value = 0.3992303
$ws.Range("AB137").Value = 0.3912847
$ws.Range("AD137").Value = 0.4385677
$ws.Range("AE137").Value = 0.3511533
$ws.Range("AF137").Value = 0.3888566
$ws.Range("AG137").Value = 0.3628673
$ws.Range("AH137").Value = 0.4464328
$ws.Range("AI137").Value = 0.2108449
$ws.Range("AJ137").Value = 0.3662842
$ws.Range("AK137").Value = 0.4051526
$ws.Range("AL137").Value = 0.422492
$ws.Range("AM137").Value = 0.3783611
$ws.Range("AN137").Value = 0.3796884
$ws.Range("AO137").Value = 0.5958856
$ws.Range("AP137").Value = 0.3795114
$ws.Range("AQ137").Value = 0.2862568
$ws.Range("AS137").Value = 0.4694077
$ws.Range("AT137").Value = 0.5159280000000001
$ws.Range("AU137").Value = 0.4651905
$ws.Range("AV137").Value = 0.472898
$ws.Range("AW137").Value = 0.5592975
$ws.Range("AX137").Value = 0.3991633
$ws.Range("AY137").Value = 0.3173978
$ws.Range("BA137").Value = 0.2428318
$ws.Range("BB137").Value = 0.3551825
$ws.Range("BC137").Value = 0.2575358
$ws.Range("BD137").Value = 0.4222794
$ws.Range("BE137").Value = 0.6428068

$ws.Range("A138").Value = "16 06 2020"
$ws.Range("B138").Value = 0.4732886
$ws.Range("C138").Value = 0.5544745
$ws.Range("D138").Value = 0.5089106
$ws.Range("F138").Value = 0.6589318
$ws.Range("G138").Value = 0.355505
$ws.Range("H138").Value = 0.2838239
$ws.Range("I138").Value = 0.2613993
$ws.Range("J138").Value = 0.2413938
$ws.Range("K138").Value = 0.3302483
$ws.Range("L138").Value = 0.3887495
$ws.Range("M138").Value = 0.4240268
$ws.Range("O138").Value = 0.3287871
$ws.Range("P138").Value = 0.4294339
$ws.Range("Q138").Value = 0.5504509
$ws.Range("R138").Value = 0.3217883
$ws.Range("S138").Value = 0.5624882
$ws.Range("T138").Value = 0.3777121
$ws.Range("U138").Value = 0.4571723
$ws.Range("V138").Value = 0.5132887
$ws.Range("W138").Value = 0.2646648
$ws.Range("X138").Value = 0.3440175
$ws.Range("Y138").Value = 0.3231518
$ws.Range("Z138").Value = 0.2438268
$ws.Range("AA138").Value = 0.3545634
$ws.Range("AB138").Value = 0.3986821
$ws.Range("AD138").Value = 0.468676
$ws.Range("AE138").Value = 0.3454702
$ws.Range("AF138").Value = 0.3939319
$ws.Range("AG138").Value = 0.3259683
$ws.Range("AH138").Value = 0.4596346
$ws.Range("AI138").Value = 0.1951879
$ws.Range("AJ138").Value = 0.3543005
$ws.Range("AK138").Value = 0.3749581
$ws.Range("AL138").Value = 0.4236747
$ws.Range("AM138").Value = 0.3631976
$ws.Range("AN138").Value = 0.3459083
$ws.Range("AO138").Value = 0.5989943
$ws.Range("AP138").Value = 0.3707371
$ws.Range("AQ138").Value = 0.2742741
$ws.Range("AS138").Value = 0.4379245
$ws.Range("AT138").Value = 0.534052
$ws.Range("AU138").Value = 0.4819388
$ws.Range("AV138").Value = 0.4956076
$ws.Range("AW138").Value = 0.5260745
$ws.Range("AX138").Value = 0.4556505
$ws.Range("AY138").Value = 0.3013279
$ws.Range("BA138").Value = 0.2555166
$ws.Range("BB138").Value = 0.374153
$ws.Range("BC138").Value = 0.2458547
$ws.Range("BD138").Value = 0.3709702
$ws.Range("BE138").Value = 0.7903167

$ws.Range("A139").Value = "17 06 2020"
$ws.Range("B139").Value = 0.4701372
$ws.Range("C139").Value = 0.5674541
$ws.Range("D139").Value = 0.4568984
$ws.Range("F139").Value = 0.6717857
$ws.Range("G139").Value = 0.3408396
$ws.Range("H139").Value = 0.2848779
$ws.Range("I139").Value = 0.2858059
$ws.Range("J139").Value = 0.0518888
$ws.Range("K139").Value = 0.3515038
$ws.Range("L139").Value = 0.3843177
$ws.Range("M139").Value = 0.4398844
$ws.Range("O139").Value = 0.3077302
$ws.Range("P139").Value = 0.4787546
$ws.Range("Q139").Value = 0.5558761
$ws.Range("R139").Value = 0.3235494
$ws.Range("S139").Value = 0.5428396
$ws.Range("T139").Value = 0.3102085
$ws.Range("U139").Value = 0.4255034
$ws.Range("V139").Value = 0.5313001000000001
$ws.Range("W139").Value = 0.2323535
$ws.Range("X139").Value = 0.328672
$ws.Range("Y139").Value = 0.3628744
$ws.Range("Z139").Value = 0.2333119
$ws.Range("AA139").Value = 0.3401903
$ws.Range("AB139").Value = 0.3758055
$ws.Range("AD139").Value = 0.4691514
$ws.Range("AE139").Value = 0.3662106
$ws.Range("AF139").Value = 0.399587
$ws.Range("AG139").Value = 0.3559375
$ws.Range("AH139").Value = 0.4014919
$ws.Range("AI139").Value = 0.1990995
$ws.Range("AJ139").Value = 0.4013537
$ws.Range("AK139").Value = 0.4443804
$ws.Range("AL139").Value = 0.4670868
$ws.Range("AM139").Value = 0.3416148
$ws.Range("AN139").Value = 0.3639225
$ws.Range("AO139").Value = 0.5542882
$ws.Range("AP139").Value = 0.3384801
$ws.Range("AQ139").Value = 0.2894988
$ws.Range("AS139").Value = 0.4460346
$ws.Range("AT139").Value = 0.5353824
$ws.Range("AU139").Value = 0.4647673
$ws.Range("AV139").Value = 0.4888657
$ws.Range("AW139").Value = 0.5547682
$ws.Range("AX139").Value = 0.4201896
$ws.Range("AY139").Value = 0.3001141
$ws.Range("BA139").Value = 0.2959499
$ws.Range("BB139").Value = 0.3649437
$ws.Range("BC139").Value = 0.2497068
$ws.Range("BD139").Value = 0.332628
$ws.Range("BE139").Value = 0.6749713000000001

$ws.Range("A140").Value = "18 06 2020"
$ws.Range("B140").Value = 0.4693962
$ws.Range("C140").Value = 0.5937305
$ws.Range("D140").Value = 0.4765241
$ws.Range("F140").Value = 0.6525919
$ws.Range("G140").Value = 0.3530451
$ws.Range("H140").Value = 0.303213
$ws.Range("I140").Value = 0.2131504
$ws.Range("J140").Value = 0.093985
$ws.Range("K140").Value = 0.4157774
$ws.Range("L140").Value = 0.3915568
$ws.Range("M140").Value = 0.4275139
$ws.Range("O140").Value = 0.3041543
$ws.Range("P140").Value = 0.5459044
$ws.Range("Q140").Value = 0.5461182
$ws.Range("R140").Value = 0.3347969
$ws.Range("S140").Value = 0.5309327
$ws.Range("T140").Value = 0.3596751
$ws.Range("U140").Value = 0.5331119
$ws.Range("V140").Value = 0.4819868
$ws.Range("W140").Value = 0.2052207
$ws.Range("X140").Value = 0.3423725
$ws.Range("Y140").Value = 0.3101311
$ws.Range("Z140").Value = 0.2528258
$ws.Range("AA140").Value = 0.3217951
$ws.Range("AB140").Value = 0.3888063
$ws.Range("AD140").Value = 0.4691002
$ws.Range("AE140").Value = 0.3480277
$ws.Range("AF140").Value = 0.3815761
$ws.Range("AG140").Value = 0.4270131
$ws.Range("AH140").Value = 0.3719538
$ws.Range("AI140").Value = 0.1869149
$ws.Range("AJ140").Value = 0.3356903
$ws.Range("AK140").Value = 0.4438576
$ws.Range("AL140").Value = 0.5279925
$ws.Range("AM140").Value = 0.360876
$ws.Range("AN140").Value = 0.3412041
$ws.Range("AO140").Value = 0.5845042
$ws.Range("AP140").Value = 0.3377094
$ws.Range("AQ140").Value = 0.2825728
$ws.Range("AS140").Value = 0.4711661
$ws.Range("AT140").Value = 0.5391277
$ws.Range("AU140").Value = 0.4841538
$ws.Range("AV140").Value = 0.4977489
$ws.Range("AW140").Value = 0.6067787
$ws.Range("AX140").Value = 0.4194234
$ws.Range("AY140").Value = 0.3310208
$ws.Range("BA140").Value = 0.3887556
$ws.Range("BB140").Value = 0.3608996
$ws.Range("BC140").Value = 0.2621056
$ws.Range("BD140").Value = 0.389045
$ws.Range("BE140").Value = 0.6101239000000001

$ws.Range("A141").Value = "19 06 2020"
$ws.Range("B141").Value = 0.4502472
$ws.Range("C141").Value = 0.627809
$ws.Range("D141").Value = 0.5400712
$ws.Range("F141").Value = 0.6997814
$ws.Range("G141").Value = 0.3413756
$ws.Range("H141").Value = 0.2995005
$ws.Range("I141").Value = 0.22768
$ws.Range("J141").Value = 0.2226463
$ws.Range("K141").Value = 0.342621
$ws.Range("L141").Value = 0.4497954
$ws.Range("M141").Value = 0.4421135
$ws.Range("O141").Value = 0.3422542
$ws.Range("P141").Value = 0.5765763
$ws.Range("Q141").Value = 0.5500533
$ws.Range("R141").Value = 0.3413952
$ws.Range("S141").Value = 0.530601
$ws.Range("T141").Value = 0.4626593
$ws.Range("U141").Value = 0.5400486
$ws.Range("V141").Value = 0.5003126
$ws.Range("W141").Value = 0.220758
$ws.Range("X141").Value = 0.2920612
$ws.Range("Y141").Value = 0.314421
$ws.Range("Z141").Value = 0.2599472
$ws.Range("AA141").Value = 0.3217902
$ws.Range("AB141").Value = 0.3442856
$ws.Range("AD141").Value = 0.4972253
$ws.Range("AE141").Value = 0.4229359
$ws.Range("AF141").Value = 0.3589104
$ws.Range("AG141").Value = 0.3713942
$ws.Range("AH141").Value = 0.4201944
$ws.Range("AI141").Value = 0.199863
$ws.Range("AJ141").Value = 0.338574
$ws.Range("AK141").Value = 0.3589787
$ws.Range("AL141").Value = 0.5108065000000001
$ws.Range("AM141").Value = 0.342528
$ws.Range("AN141").Value = 0.3297343
$ws.Range("AO141").Value = 0.5267998
$ws.Range("AP141").Value = 0.291829
$ws.Range("AQ141").Value = 0.2947598
$ws.Range("AS141").Value = 0.4815724
$ws.Range("AT141").Value = 0.6089183
$ws.Range("AU141").Value = 0.5320142
$ws.Range("AV141").Value = 0.5144501
$ws.Range("AW141").Value = 0.6255086
$ws.Range("AX141").Value = 0.4517852
$ws.Range("AY141").Value = 0.3534276
$ws.Range("BA141").Value = 0.4111319
$ws.Range("BB141").Value = 0.3513644
$ws.Range("BC141").Value = 0.3239032
$ws.Range("BD141").Value = 0.425462
$ws.Range("BE141").Value = 0.4570312

$ws.Range("A142").Value = "20 06 2020"
$ws.Range("B142").Value = 0.3584229
$ws.Range("C142").Value = 0.6410641
$ws.Range("D142").Value = 0.5962853
$ws.Range("F142").Value = 0.7922577
$ws.Range("G142").Value = 0.3658923
$ws.Range("H142").Value = 0.2725695
$ws.Range("I142").Value = 0.2208733
$ws.Range("J142").Value = 0.2286446
$ws.Range("K142").Value = 0.3166596
$ws.Range("L142").Value = 0.4798122
$ws.Range("M142").Value = 0.4730118
$ws.Range("O142").Value = 0.3884514
$ws.Range("P142").Value = 0.6496303
$ws.Range("Q142").Value = 0.6020694
$ws.Range("R142").Value = 0.3118412
$ws.Range("S142").Value = 0.5268617
$ws.Range("T142").Value = 0.5149603
$ws.Range("U142").Value = 0.5203803
$ws.Range("V142").Value = 0.511292
$ws.Range("W142").Value = 0.2270082
$ws.Range("X142").Value = 0.2818457
$ws.Range("Y142").Value = 0.3056641
$ws.Range("Z142").Value = 0.2479496
$ws.Range("AA142").Value = 0.3195018
$ws.Range("AB142").Value = 0.3215328
$ws.Range("AD142").Value = 0.6292158
$ws.Range("AE142").Value = 0.4601154
$ws.Range("AF142").Value = 0.3932228
$ws.Range("AG142").Value = 0.4004946
$ws.Range("AH142").Value = 0.4175432
$ws.Range("AI142").Value = 0.2033112
$ws.Range("AJ142").Value = 0.3446759
$ws.Range("AK142").Value = 0.2975971
$ws.Range("AL142").Value = 0.626241
$ws.Range("AM142").Value = 0.321261
$ws.Range("AN142").Value = 0.3199374
$ws.Range("AO142").Value = 0.5104937000000001
$ws.Range("AP142").Value = 0.2871372
$ws.Range("AQ142").Value = 0.2981611
$ws.Range("AS142").Value = 0.4665762
$ws.Range("AT142").Value = 0.6570602
$ws.Range("AU142").Value = 0.5137547
$ws.Range("AV142").Value = 0.518134
$ws.Range("AW142").Value = 0.6770944
$ws.Range("AX142").Value = 0.4228395
$ws.Range("AY142").Value = 0.3620121
$ws.Range("BA142").Value = 0.3822627
$ws.Range("BB142").Value = 0.3061836
$ws.Range("BC142").Value = 0.2982414
$ws.Range("BD142").Value = 0.3298753
$ws.Range("BE142").Value = 0.5293353

$ws.Range("A143").Value = "21 06 2020"

$ws.Range("A144").Value = "22 06 2020"
